$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(3,6,10,12,13,14,15,18,19,20,21,22,24,30,33,37,39,40,41,42,45,46,47,48,49,51,57,60,64,66,67,68,69,72,73,74,75,76,78,86,87,88,89,93,95,102,112,113,114,115,119,121,128,138,139,140,141,145,147,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
